$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("N37").Copy()
$ws.Range("AJ33").PasteSpecial(-4122)
